$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "62.889.61"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +5.06%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.437.00"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +5.77%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.07%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "565.58"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +4.71%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "140.02"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +8.12%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.03%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.588"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +2.77%  "
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.435.88"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +5.79%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.105"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +4.33%  "
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +4.84%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.151"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +0.53%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.350"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +5.87%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.28"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  +13.02%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.872.98"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +5.93%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "62.756.02"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  +4.88%  "
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  +7.93%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.441.50"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +6.27%  "
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +7.74%  "
$c.ClearFormats()

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "339.44"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +8.91%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +4.19%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.85"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +4.90%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.57"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -1.91%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "65.61"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  +3.70%  "
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +2.45%  "
$c.ClearFormats()

$c = $ws.Range("B27")
$c.NumberFormat = "@"
$c.Value = "Binance-PegBSC-USD"
$c.ClearFormats()
$c = $ws.Range("C27")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +0.14%  "
$c.ClearFormats()

$c = $ws.Range("B28")
$c.NumberFormat = "@"
$c.Value = "Fetch.AI"
$c.ClearFormats()
$c = $ws.Range("C28")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c.ClearFormats()
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.53"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  +13.75%  "
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.22"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +6.77%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.37"
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  +16.70%  "
$c.ClearFormats()

$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = "PEPE"
$c.ClearFormats()
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c.ClearFormats()
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0790"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  +9.31%  "
$c.ClearFormats()

$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = "PancakeSwap"
$c.ClearFormats()
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c.ClearFormats()
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.82"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +6.23%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.52"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +12.00%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "172.98"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +0.86%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.47"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +8.95%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.397"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +5.41%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "382.49"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +20.78%  "
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.63"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +5.47%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  +12.29%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -0.03%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +13.17%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "40.01"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +5.59%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "144.92"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +6.76%  "
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +7.55%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "20.60"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +10.34%  "
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.595"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +5.07%  "
$c.ClearFormats()

$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = "Stellar"
$c.ClearFormats()
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c.ClearFormats()
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0954"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +1.81%  "
$c.ClearFormats()

$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = "Hedera"
$c.ClearFormats()
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c.ClearFormats()
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0520"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +6.69%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0224"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +5.23%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.88"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +7.19%  "
$c.ClearFormats()
